$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Perfect 100" column header in F1 (new shared string)
$ws.Range("F1").Value = "Perfect 100"
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").VerticalAlignment = -4108

# Add F4/F5/F6/F8 values of 100 (new "Perfect 100" data column)
$ws.Range("F4").Value = 100
$ws.Range("F5").Value = 100
$ws.Range("F6").Value = 100
$ws.Range("F8").Value = 100

# Adjust QPSK 2400 integral limit row (row 8): baud 2400 -> 1200, integral limit 30 -> 32
$ws.Range("A8").Value = 1200
$ws.Range("D8").Value = 32

# Update selected cell to D9
$ws.Range("D9").Select()
